$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.973.16"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.769.47"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "328.89"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  -0.03%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4664"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.98%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3518"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.92%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "44.01"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.14%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07389"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.083"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.94%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  -1.02%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.015"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.197"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "1.767.37"
$ws.Range("E16").Value = "  -0.53%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "92.27"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("E18").Value = "  -0.37%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06423"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.03%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "16.92"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.794"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "28.012.42"
$ws.Range("E23").Value = "  +0.79%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.14"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.59%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.156"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.54%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "163.60"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").Value = "1.967.66"
$ws.Range("E28").Value = "  -0.45%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.206"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.61%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "123.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.96%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.076"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.49%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09322"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.23%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.654"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.40%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.550"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  -1.69%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02265"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.34%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06112"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2070"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.909"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.191"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6156"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.436"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.37%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "7.766"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.02%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "13.19"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.37%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.743"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5801"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.97%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "123.90"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -1.13%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06814"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "72.04"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
